$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a small 2-row x 4-column dataframe into the sheet starting at B3
# (row 3, col B), i.e. Range("B3:E4"). Its first row overlaps/updates the
# existing row 3 data (and extends it one column to the right, into E3);
# its second row lands on the previously-empty row 4.
#
#        B        C        D        E
#  3     0        2        6        3
#  4     z        0        2        6
#
# All values in this sheet are stored as text (matching the surrounding
# cells, e.g. "1", "2", "3" are text, not numbers), so numeric-looking
# strings need NumberFormat "@" (Text) set first - otherwise Excel's COM
# layer auto-converts a typed "0"/"2"/"6"/"3" into a real number.

$dfRange = $ws.Range("B3:E4")
$dfRange.NumberFormat = "@"

$dfValues = New-Object 'object[,]' 2,4
$dfValues[0,0] = "0"; $dfValues[0,1] = "2"; $dfValues[0,2] = "6"; $dfValues[0,3] = "3"
$dfValues[1,0] = "z"; $dfValues[1,1] = "0"; $dfValues[1,2] = "2"; $dfValues[1,3] = "6"

$dfRange.Value = $dfValues

# One extra cell belonging to the same write, just above the block: E2.
$e2 = $ws.Cells.Item(2, 5)
$e2.NumberFormat = "@"
$e2.Value = "3"
